$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.533.84"
$ws.Range("D3").Value = "3.495.71"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "612.54"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "187.13"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.213"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.05"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000308"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.51"
$ws.Range("D14").Value = "4.052.83"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "601.93"
$ws.Range("D16").Value = "69.580.35"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.91"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.62"
$ws.Range("D19").Value = "3.497.60"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.987"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.22"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "105.68"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.14"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.65"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.94"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.73"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.48"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.19"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.44"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.45"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.18"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "517.40"
$ws.Range("D40").Value = "3.609.40"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.82"
$ws.Range("D42").Value = "0.0₃0778"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.94"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.33"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.78"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.38"

$ws.Range("E2").Value = "  -1.84%  "
$ws.Range("E3").Value = "  -2.08%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("E5").Value = "  +5.24%  "
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  -5.25%  "
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("E11").Value = "  -2.96%  "
$ws.Range("E12").Value = "  -3.93%  "
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("E14").Value = "  -2.04%  "
$ws.Range("E15").Value = "  +4.94%  "
$ws.Range("E16").Value = "  -1.83%  "
$ws.Range("E17").Value = "  -1.66%  "
$ws.Range("E18").Value = "  -1.69%  "
$ws.Range("E19").Value = "  -2.00%  "
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("E21").Value = "  -1.58%  "
$ws.Range("E22").Value = "  -2.51%  "
$ws.Range("E23").Value = "  +12.67%  "
$ws.Range("E24").Value = "  +4.73%  "
$ws.Range("E25").Value = "  +1.86%  "
$ws.Range("E26").Value = "  +3.04%  "
$ws.Range("E27").Value = "  -2.60%  "
$ws.Range("E28").Value = "  +5.55%  "
$ws.Range("E29").Value = "  +3.12%  "
$ws.Range("E30").Value = "  -3.53%  "
$ws.Range("E31").Value = "  +15.78%  "
$ws.Range("E32").Value = "  +1.05%  "
$ws.Range("E33").Value = "  -1.31%  "
$ws.Range("E34").Value = "  +0.57%  "
$ws.Range("E35").Value = "  -6.49%  "
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("E37").Value = "  +7.74%  "
$ws.Range("E38").Value = "  -5.14%  "
$ws.Range("E39").Value = "  -4.73%  "
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("E41").Value = "  -3.69%  "
$ws.Range("E42").Value = "  -3.76%  "
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("E44").Value = "  -0.98%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("E46").Value = "  +2.92%  "
$ws.Range("E48").Value = "  -6.16%  "
$ws.Range("E49").Value = "  +0.32%  "
$ws.Range("E50").Value = "  -2.75%  "
$ws.Range("E51").Value = "  -8.72%  "
